$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @(3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,3,1,2.991953333333333,8.97586,0.3033492052285701,0.3033492052285701,227.6353692687556,2048.7183234188,0.2908456054291093,0.2908456054291093),
    @(3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,3,1,4.959157333333334,14.877472,0.5028007686183057,0.5028007686183057,377.3052200575289,3395.74698051776,0.4820760741694525,0.4820760741694525),
    @(3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,3,1,1.911955666666667,5.735867,0.1938500261531243,0.1938500261531243,145.4664179946511,1309.19776195186,0.1858598184770984,0.1858598184770984),
    @(3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,3,1,2.991953333333333,8.97586,0.3033492052285701,0.3033492052285701,2.796701450753333,25.17031305678,0.003573294999198823,0.003573294999198822),
    @(3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,3,1,4.959157333333334,14.877472,0.5028007686183057,0.5028007686183057,4.635527684917333,41.719749164256,0.005922730111467926,0.005922730111467926),
    @(3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,3,1,1.911955666666667,5.735867,0.1938500261531243,0.1938500261531243,1.787183351815667,16.084650166341,0.002283451932981302,0.002283451932981302),
    @(3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,3,1,2.991953333333333,8.97586,0.3033492052285701,0.3033492052285701,6.989458300017779,62.90512470016001,0.008930304800262013,0.008930304800262011),
    @(3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,3,1,4.959157333333334,14.877472,0.5028007686183057,0.5028007686183057,11.58501471209245,104.265132408832,0.01480196433738535,0.01480196433738535),
    @(3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,3,1,1.911955666666667,5.735867,0.1938500261531243,0.1938500261531243,4.466491590883556,40.198424317952,0.005706755743044618,0.005706755743044617)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}

Write-Output "applied NATMI recompute updates to E2:T10"
